$wb = $excel.ActiveWorkbook

function Update-EventSheet($ws) {
    # --- Step 1: insert the rows needed so the table grows from 8 data rows to 11 ---
    # Insert two new blank rows before the current row 8 (银魂only / row A=7)
    $ws.Rows.Item(8).Insert()
    $ws.Rows.Item(8).Insert()
    # Insert one more blank row before the current row 11 (梦时空SPO1动漫展, now shifted down)
    $ws.Rows.Item(11).Insert()

    # Fix up the style of column A for the freshly inserted rows so it matches
    # the bordered/bold/centered style used throughout column A (copy format from A1)
    $ws.Range("A1").Copy()
    $ws.Range("A8").PasteSpecial(-4122)
    $ws.Range("A9").PasteSpecial(-4122)
    $ws.Range("A11").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # --- Step 2: update the "want-to-go count" numbers that changed on existing rows ---
    $ws.Range("F3").Value = 2570
    $ws.Range("F4").Value = 516
    $ws.Range("F6").Value = 6566
    $ws.Range("F7").Value = 400

    # --- Step 3: fill in row 8 (new entry: 柯暮卿 inner venue) ---
    $ws.Range("A8").Value = 7
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = "2024.04.04"
    $ws.Range("C8").Value = "合肥·第二届漫画城市动漫展内场-柯暮卿"
    $ws.Range("D8").Value = "凤淮路与固镇路西北角 庐阳全民健身中心"
    $ws.Range("E8").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("F8").Value = 4
    $ws.Range("G8").Value = 158
    $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=82192"
    $ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg"

    # --- Step 4: fill in row 9 (new entry: 风袖 inner venue) ---
    $ws.Range("A9").Value = 8
    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = "2024.04.04"
    $ws.Range("C9").Value = "合肥·第二届漫画城市动漫展内场-风袖"
    $ws.Range("D9").Value = "凤淮路与固镇路西北角 庐阳全民健身中心"
    $ws.Range("E9").Value = "2024.04.04 10:00-04.04 17:00"
    $ws.Range("F9").Value = 5
    $ws.Range("G9").Value = 158
    $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=82191"
    $ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg"

    # --- Step 5: row 10 already holds the old row-8 data (银魂only); just update A + F ---
    $ws.Range("A10").Value = 9
    $ws.Range("F10").Value = 7

    # --- Step 6: fill in row 11 (new entry: BH动漫游戏展) ---
    $ws.Range("A11").Value = 10
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2024.05.03"
    $ws.Range("C11").Value = "合肥·BH动漫游戏展"
    $ws.Range("D11").Value = "科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地"
    $ws.Range("E11").Value = "2024.05.03 10:00-05.04 16:00"
    $ws.Range("F11").Value = 1
    $ws.Range("G11").Value = 40
    $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=82199"
    $ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg"

    # --- Step 7: row 12 already holds the old row-9 data (梦时空SPO1动漫展); update A + F ---
    $ws.Range("A12").Value = 11
    $ws.Range("F12").Value = 134
}

$wsExhibition = $wb.Worksheets.Item("展览")
Update-EventSheet $wsExhibition

$wsAllTypes = $wb.Worksheets.Item("全部类型")
Update-EventSheet $wsAllTypes
